$wb = $excel.ActiveWorkbook
try {
  $wb.ApplyTheme("Office")
  Write-Host "wb ApplyTheme ok"
} catch {
  Write-Host "ERR1: $_"
}
